$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# The edit moves the lone "_GoBack" bookmark from the end of the document
# (right after the superscript "." that closes the citation/back-reference
# paragraph) to a spot in the middle of the "objective of this study" run,
# splitting that run right after "...determining the us|e of offensive...".
#
# Word keeps at most one bookmark with a given name, so re-adding a
# bookmark named "_GoBack" at the new location automatically removes it
# from its old location - exactly matching both hunks of the diff (the
# insertion in the first paragraph and the deletion near "Stanford
# Parser ... [11]." later in the document).
# ---------------------------------------------------------------------------

$rng = $d.Content
$found = $rng.Find.Execute("finding and determining the us", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)

if ($found) {
    $splitPoint = $rng.End
    $bmRange = $d.Range($splitPoint, $splitPoint)
    $d.Bookmarks.Add("_GoBack", $bmRange)
    Write-Host "Moved _GoBack bookmark to offset $splitPoint"
} else {
    Write-Host "Could not locate split point for _GoBack bookmark"
}
